$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A - case identifiers
$ws.Range("A26").Value = "Caso #20"
$ws.Range("A27").Value = "Caso #21"
$ws.Range("A28").Value = "Caso #22"

# Column B - action names
$ws.Range("B26").Value = "Ingresar compra"
$ws.Range("B27").Value = "Ver compras usuario"
$ws.Range("B28").Value = "Ver compras admin"

# Column C - expected results
$ws.Range("C26").Value = "Ingrersar compra desde menu de usuario"
$ws.Range("C27").Value = "Ver compras hechas por el usuario logueado"
$ws.Range("C28").Value = "Ver las compras de todos los usuarios"

# Column D - dates
$ws.Range("D26").Value = 45274
$ws.Range("D27").Value = 45274
$ws.Range("D28").Value = 45274

# Column E - approved?
$ws.Range("E26").Value = "SI"
$ws.Range("E27").Value = "SI"
$ws.Range("E28").Value = "SI"

# Column F - observations
$ws.Range("F26").Value = "-"
$ws.Range("F27").Value = "-"
$ws.Range("F28").Value = "-"

# Column G - tested result
$ws.Range("G26").Value = "OK"
$ws.Range("G27").Value = "OK"
$ws.Range("G28").Value = "OK"

# Update selection to match the new active range
$ws.Range("E26:G28").Select()

$wb.Save()
